# Natmi following Dr Hou advice:
# Add the "ECs" cluster into the Rarres2-Ccrl2 sending/target cluster grid,
# expanding the LR-pair table from a 2x2 (FAPs/sCs) combination to the full
# 3x3 (ECs/FAPs/sCs) sending x target cluster combination.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:T10").ClearContents()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rarres2"
$ws.Range("C2").Value = "Ccrl2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.445726333333333
$ws.Range("H2").Value = 4.337179
$ws.Range("I2").Value = 0.01544456920939864
$ws.Range("J2").Value = 0.01544456920939864
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 179.3411446666667
$ws.Range("N2").Value = 538.0234340000001
$ws.Range("O2").Value = 0.9868930934198245
$ws.Range("P2").Value = 0.9868930934198243
$ws.Range("Q2").Value = 259.2782154947429
$ws.Range("R2").Value = 2333.503939452686
$ws.Range("S2").Value = 0.01524213868359999
$ws.Range("T2").Value = 0.01524213868359999

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rarres2"
$ws.Range("C3").Value = "Ccrl2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.445726333333333
$ws.Range("H3").Value = 4.337179
$ws.Range("I3").Value = 0.01544456920939864
$ws.Range("J3").Value = 0.01544456920939864
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.5124153333333333
$ws.Range("N3").Value = 1.537246
$ws.Range("O3").Value = 0.002819760933103243
$ws.Range("P3").Value = 0.002819760933103243
$ws.Range("Q3").Value = 0.7408123410037777
$ws.Range("R3").Value = 6.667311069034001
$ws.Range("S3").Value = 0.000043549992885271517377085426
$ws.Range("T3").Value = 0.000043549992885271517377085426

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rarres2"
$ws.Range("C4").Value = "Ccrl2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.445726333333333
$ws.Range("H4").Value = 4.337179
$ws.Range("I4").Value = 0.01544456920939864
$ws.Range("J4").Value = 0.01544456920939864
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.869410666666667
$ws.Range("N4").Value = 5.608232
$ws.Range("O4").Value = 0.01028714564707241
$ws.Range("P4").Value = 0.01028714564707241
$ws.Range("Q4").Value = 2.702656228614222
$ws.Range("R4").Value = 24.323906057528
$ws.Range("S4").Value = 0.0001588805329133737
$ws.Range("T4").Value = 0.0001588805329133737

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rarres2"
$ws.Range("C5").Value = "Ccrl2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 75.55280566666666
$ws.Range("H5").Value = 226.658417
$ws.Range("I5").Value = 0.8071240795570661
$ws.Range("J5").Value = 0.8071240795570661
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 179.3411446666667
$ws.Range("N5").Value = 538.0234340000001
$ws.Range("O5").Value = 0.9868930934198245
$ws.Range("P5").Value = 0.9868930934198243
$ws.Range("Q5").Value = 13549.72665103822
$ws.Range("R5").Value = 121947.539859344
$ws.Range("S5").Value = 0.7965451796477014
$ws.Range("T5").Value = 0.7965451796477013

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rarres2"
$ws.Range("C6").Value = "Ccrl2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 75.55280566666666
$ws.Range("H6").Value = 226.658417
$ws.Range("I6").Value = 0.8071240795570661
$ws.Range("J6").Value = 0.8071240795570661
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.5124153333333333
$ws.Range("N6").Value = 1.537246
$ws.Range("O6").Value = 0.002819760933103243
$ws.Range("P6").Value = 0.002819760933103243
$ws.Range("Q6").Value = 38.71441609995355
$ws.Range("R6").Value = 348.429744899582
$ws.Range("S6").Value = 0.002275896947701929
$ws.Range("T6").Value = 0.002275896947701929

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rarres2"
$ws.Range("C7").Value = "Ccrl2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 75.55280566666666
$ws.Range("H7").Value = 226.658417
$ws.Range("I7").Value = 0.8071240795570661
$ws.Range("J7").Value = 0.8071240795570661
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.869410666666667
$ws.Range("N7").Value = 5.608232
$ws.Range("O7").Value = 0.01028714564707241
$ws.Range("P7").Value = 0.01028714564707241
$ws.Range("Q7").Value = 141.2392208098604
$ws.Range("R7").Value = 1271.152987288744
$ws.Range("S7").Value = 0.008303002961662797
$ws.Range("T7").Value = 0.008303002961662795

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Rarres2"
$ws.Range("C8").Value = "Ccrl2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 16.60889166666666
$ws.Range("H8").Value = 49.82667499999999
$ws.Range("I8").Value = 0.1774313512335352
$ws.Range("J8").Value = 0.1774313512335352
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 179.3411446666667
$ws.Range("N8").Value = 538.0234340000001
$ws.Range("O8").Value = 0.9868930934198245
$ws.Range("P8").Value = 0.9868930934198243
$ws.Range("Q8").Value = 2978.657643144661
$ws.Range("R8").Value = 26807.91878830195
$ws.Range("S8").Value = 0.1751057750885229
$ws.Range("T8").Value = 0.1751057750885229

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Rarres2"
$ws.Range("C9").Value = "Ccrl2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 16.60889166666666
$ws.Range("H9").Value = 49.82667499999999
$ws.Range("I9").Value = 0.1774313512335352
$ws.Range("J9").Value = 0.1774313512335352
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.5124153333333333
$ws.Range("N9").Value = 1.537246
$ws.Range("O9").Value = 0.002819760933103243
$ws.Range("P9").Value = 0.002819760933103243
$ws.Range("Q9").Value = 8.510650759672222
$ws.Range("R9").Value = 76.59585683705
$ws.Range("S9").Value = 0.0005003139925160424
$ws.Range("T9").Value = 0.0005003139925160424

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Rarres2"
$ws.Range("C10").Value = "Ccrl2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 16.60889166666666
$ws.Range("H10").Value = 49.82667499999999
$ws.Range("I10").Value = 0.1774313512335352
$ws.Range("J10").Value = 0.1774313512335352
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.869410666666667
$ws.Range("N10").Value = 5.608232
$ws.Range("O10").Value = 0.01028714564707241
$ws.Range("P10").Value = 0.01028714564707241
$ws.Range("Q10").Value = 31.04883924317778
$ws.Range("R10").Value = 279.4395531886
$ws.Range("S10").Value = 0.001825262152496237
$ws.Range("T10").Value = 0.001825262152496236

